# Generate Report for Handoff
# Updates the localization-status workbook with a freshly generated handoff
# report: new source-file GUID, new xliff checksum/GUID and refreshed
# timestamps for the handoff that just completed.

$wb = $excel.ActiveWorkbook

$oldGuid = "27185796-f654-467d-86dc-5a1a446f2d03"
$newGuid = "68aeb95c-19ca-4db1-a4f6-7c5a94a4946d"

$oldXliffHash = "df967b91c0884aaa958ea51165d2b53c54d28072"
$newXliffHash = "380230b49e7c96d04f0098e0e5e9e7a06b1b4e47"

$oldGenDate = "2016-08-28 12:58:10"
$newGenDate = "2016-08-28 12:58:36"

$oldZhHandoffDate = "2016-08-28 12:58:05"
$newZhHandoffDate = "2016-08-28 12:58:32"

# The hyperlink target (commit-pinned GitHub URL) is not part of this change,
# only the visible display text / file names move to the new GUID.
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/388dacd3b852fe5145a7512f9851adff6d16f5dc/e2e/$oldGuid.md"

function Set-HyperlinkDisplay {
    param($ws, $rangeAddr, $displayText)

    $r = $ws.Range($rangeAddr)
    $r.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($r, $hyperlinkAddress, "", "", $displayText) | Out-Null
}

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
Set-HyperlinkDisplay $wsOverview "B2" "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = $newGenDate

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
Set-HyperlinkDisplay $wsZh "A2" "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newXliffHash.zh-cn.xlf"
$wsZh.Range("H2").Value = $newZhHandoffDate

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
Set-HyperlinkDisplay $wsDe "A2" "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newXliffHash.de-de.xlf"
$wsDe.Range("H2").Value = $newGenDate
